$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.110.64'
$ws.Range('E2').Value = '  -2.36%  '

$ws.Range('D3').Value = '1.898.40'
$ws.Range('E3').Value = '  -2.72%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '332.32'
$ws.Range('E5').Value = '  -2.96%  '

$ws.Range('E6').Value = '  -0.04%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4595'
$ws.Range('E7').Value = '  -3.67%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4129'
$ws.Range('E8').Value = '  -0.25%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.87'
$ws.Range('E9').Value = '  +0.04%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08003'
$ws.Range('E10').Value = '  -3.04%  '

$ws.Range('E11').Value = '  -2.44%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.20'
$ws.Range('E12').Value = '  -2.14%  '

$ws.Range('D13').Value = '1.892.03'
$ws.Range('E13').Value = '  -2.99%  '

$ws.Range('E14').Value = '  -3.96%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.119'
$ws.Range('E15').Value = '  -4.01%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  -0.05%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '89.08'
$ws.Range('E17').Value = '  -3.35%  '

$ws.Range('E18').Value = '  -2.98%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06576'
$ws.Range('E19').Value = '  -1.74%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.63'
$ws.Range('E20').Value = '  -2.41%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.21%  '

$ws.Range('D22').Value = '29.059.99'
$ws.Range('E22').Value = '  -2.42%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.492'
$ws.Range('E23').Value = '  -1.67%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.41'
$ws.Range('E24').Value = '  +1.23%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.199'
$ws.Range('E25').Value = '  -2.66%  '

$ws.Range('D26').Value = '2.121.24'

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '156.61'
$ws.Range('E27').Value = '  -3.03%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.76'
$ws.Range('E28').Value = '  -2.21%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.119'
$ws.Range('E29').Value = '  -2.99%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.628'
$ws.Range('E30').Value = '  -1.33%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '116.98'
$ws.Range('E31').Value = '  -4.72%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.050'
$ws.Range('E32').Value = '  +3.73%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09386'
$ws.Range('E33').Value = '  -2.42%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.420'
$ws.Range('E34').Value = '  -3.89%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.540'
$ws.Range('E35').Value = '  -4.03%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.350'
$ws.Range('E36').Value = '  -3.18%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06090'
$ws.Range('E37').Value = '  -3.53%  '

$ws.Range('E38').Value = '  -3.59%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.425'
$ws.Range('E39').Value = '  -1.13%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.177'
$ws.Range('E40').Value = '  -0.72%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5839'
$ws.Range('E41').Value = '  -4.07%  '

$ws.Range('E42').Value = '  -0.01%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1828'
$ws.Range('E43').Value = '  -3.32%  '

$ws.Range('E44').Value = '  -5.43%  '

$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.314'
$ws.Range('E45').Value = '  -2.18%  '

$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.233'
$ws.Range('E46').Value = '  -2.71%  '

$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.07511'
$ws.Range('E47').Value = '  +2.38%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '12.15'
$ws.Range('E48').Value = '  -2.76%  '

$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.5535'
$ws.Range('E49').Value = '  -3.12%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.925'
$ws.Range('E50').Value = '  -3.19%  '

$ws.Range('E51').Value = '  -1.54%  '
